# Update the client's phone number, address and birth date shown on the
# form (both the horizontal row-4 layout and the vertical A/B layout
# mirror the same fields):
#   phone number   "586"              -> "261"
#   address        "Нарын"            -> "Бишкек"
#   birth date     1999-10-14 (36447) -> 1996-08-09 (35286)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Birth date (stored as a date serial number) - simple numeric assignment.
$ws.Range("F4").Value = 35286
$ws.Range("B5").Value = 35286

# Address - plain text, no ambiguity with numeric auto-detection.
$ws.Range("H4").Value = "Бишкек"
$ws.Range("B7").Value = "Бишкек"

# Phone number - the new value "261" looks like a number, but the source
# file stores it as text, so a plain .Value assignment would silently turn
# it into a numeric cell. Stage it as text in a scratch cell (forced with a
# Text number format), copy it, and paste-special "values only" into the
# target cells so the destination keeps its own formatting/style while the
# cell's stored type becomes text, matching the original layout.
$scratch = $ws.Range("Z1")
$scratch.NumberFormat = "@"
$scratch.Value = "261"
$scratch.Copy()

$ws.Range("G4").PasteSpecial(-4163)
$ws.Range("B6").PasteSpecial(-4163)

$scratch.Clear()
